$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 7 first (2005, Iwate Prefecture Environmental Health Research Center),
# then row 3 (2022, Gifu Prefecture Health Environment Research Institute),
# shifting remaining rows up so the final table is A1:C6.
$ws.Rows(7).Delete(4162)
$ws.Rows(3).Delete(4162)
